$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# D-column text values that look numeric are forced to stay as text
# (matching the source inlineStr cells) by temporarily applying a text
# number format before assignment, then resetting the style so no stray
# style index is left behind.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "58.173.56"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -4.32%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.642.69"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -2.73%  "

$ws.Range("E4").Value = "  -0.11%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "520.67"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -1.12%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "143.95"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.60%  "

$ws.Range("E7").Value = "  +0.20%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.569"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -1.79%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "6.65"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -0.04%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.102"
$cell.Style = "Normal"

$ws.Range("E11").Value = "  -1.00%  "

$ws.Range("E12").Value = "  +1.51%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "3.106.52"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -1.91%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "58.197.10"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -4.24%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "20.83"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -2.32%  "

$ws.Range("E16").Value = "  -1.65%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "2.646.59"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -2.01%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "337.62"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -3.30%  "

$ws.Range("E19").Value = "  -2.66%  "

$ws.Range("E20").Value = "  -1.31%  "

$ws.Range("E21").Value = "  -0.98%  "

$ws.Range("E22").Value = "  +0.02%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "64.54"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +1.42%  "

$ws.Range("E24").Value = "  +0.67%  "

$ws.Range("E25").Value = "  -1.58%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +0.66%  "

$ws.Range("E27").Value = "  -3.05%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "7.10"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -2.89%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "6.67"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -2.26%  "

$ws.Range("E30").Value = "  +0.07%  "

$ws.Range("E31").Value = "  -1.03%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "152.47"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +1.66%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "18.81"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -1.71%  "

$ws.Range("E34").Value = "  -3.04%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.911"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -3.91%  "

$ws.Range("E36").Value = "  -4.89%  "

$ws.Range("E37").Value = "  -2.79%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "36.78"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -0.56%  "

$ws.Range("E39").Value = "  -5.92%  "

$ws.Range("E40").Value = "  -1.04%  "

$ws.Range("E41").Value = "  +0.18%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.605"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -1.00%  "

$ws.Range("E43").Value = "  -2.52%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "269.51"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -6.05%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "19.38"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -3.20%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.0538"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -0.65%  "

$ws.Range("E47").Value = "  +1.46%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "2.039.46"
$cell.Style = "Normal"

$ws.Range("E49").Value = "  -3.41%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "4.65"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -3.22%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "18.30"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -4.10%  "
